# Auto-generated Excel COM-interop script
# Applies per-cell value updates to match the target OOXML diff
# (scheduled Sheets runner: recalculated currentAveragePrice/profit columns)

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 195.6
$ws.Range("I33").Value = 243.9
$ws.Range("K33").Value = 243.9
$ws.Range("M33").Value = -14.90000000000001

$ws.Range("H41").Value = 17243738
$ws.Range("J41").Value = 41671650
$ws.Range("L41").Value = 41671650
$ws.Range("N41").Value = -41672530

$ws.Range("M64").ClearContents()
$ws.Range("H64").Value = 10500.75
$ws.Range("I64").Value = 0
$ws.Range("J64").Value = 10500.75
$ws.Range("K64").Value = 0
$ws.Range("L64").Value = 10500.75
$ws.Range("N64").Value = -10996.75

$ws.Range("M67").ClearContents()
$ws.Range("H67").Value = 10500.75
$ws.Range("I67").Value = 0
$ws.Range("J67").Value = 10500.75
$ws.Range("K67").Value = 0
$ws.Range("L67").Value = 10500.75
$ws.Range("N67").Value = -12216.75

$ws.Range("H116").Value = 3772.5454
$ws.Range("I116").Value = 3687.375
$ws.Range("K116").Value = 3687.375
$ws.Range("M116").Value = -245.375

$ws.Range("H131").Value = 3068.9333
$ws.Range("I131").Value = 2114.2
$ws.Range("J131").Value = 4978.4
$ws.Range("K131").Value = 6342.599999999999
$ws.Range("L131").Value = 14935.2
$ws.Range("M131").Value = -1302.599999999999
$ws.Range("N131").Value = -25015.2

$ws.Range("H132").Value = 4024.0815
$ws.Range("I132").Value = 4050.743
$ws.Range("K132").Value = 12152.229
$ws.Range("M132").Value = -9622.228999999999

$ws.Range("H135").Value = 3459.3076
$ws.Range("I135").Value = 1219
$ws.Range("K135").Value = 10971
$ws.Range("M135").Value = -8436

$ws.Range("H137").Value = 3879.861
$ws.Range("I137").Value = 3483.7827
$ws.Range("J137").Value = 4580.615
$ws.Range("K137").Value = 10451.3481
$ws.Range("L137").Value = 13741.845
$ws.Range("M137").Value = -7901.348100000001
$ws.Range("N137").Value = -18841.845

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 33334.72
$ws.Range("I2").Value = 41693.32
$ws.Range("J2").Value = 3482.5715
$ws.Range("K2").Value = 41693.32
$ws.Range("L2").Value = 3482.5715
$ws.Range("M2").Value = -41580.32
$ws.Range("N2").Value = -3708.5715

$ws.Range("H32").Value = 4742.326
$ws.Range("I32").Value = 4661.0444
$ws.Range("J32").Value = 8400
$ws.Range("K32").Value = 4661.0444
$ws.Range("L32").Value = 8400
$ws.Range("M32").Value = -4374.0444
$ws.Range("N32").Value = -8974

$ws.Range("H45").Value = 6818.263
$ws.Range("I45").Value = 2813.25
$ws.Range("K45").Value = 2813.25
$ws.Range("M45").Value = -2436.25

$ws.Range("H112").Value = 49500
$ws.Range("J112").Value = 49500
$ws.Range("L112").Value = 49500
$ws.Range("N112").Value = -52454

$ws.Range("H116").Value = 33334.72
$ws.Range("I116").Value = 41693.32
$ws.Range("J116").Value = 3482.5715
$ws.Range("K116").Value = 41693.32
$ws.Range("L116").Value = 3482.5715
$ws.Range("M116").Value = -39399.32
$ws.Range("N116").Value = -8070.5715

$ws.Range("H128").Value = 62750
$ws.Range("J128").Value = 62750
$ws.Range("L128").Value = 62750
$ws.Range("N128").Value = -72710

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 33334.72
$ws.Range("I3").Value = 41693.32
$ws.Range("J3").Value = 3482.5715
$ws.Range("K3").Value = 41693.32
$ws.Range("L3").Value = 3482.5715
$ws.Range("M3").Value = -41579.32
$ws.Range("N3").Value = -3710.5715

$ws.Range("H99").Value = 2218.4614
$ws.Range("I99").Value = 2034.2
$ws.Range("J99").Value = 2832.6667
$ws.Range("K99").Value = 2034.2
$ws.Range("L99").Value = 2832.6667
$ws.Range("M99").Value = -536.2
$ws.Range("N99").Value = -5828.6667

$ws.Range("N107").ClearContents()
$ws.Range("H107").Value = 1000
$ws.Range("I107").Value = 1000
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 1000
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = 920

$ws.Range("H134").Value = 2363.682
$ws.Range("I134").Value = 2285.762
$ws.Range("K134").Value = 6857.286
$ws.Range("M134").Value = -4322.286

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2754.6274
$ws.Range("J31").Value = 2783.5
$ws.Range("L31").Value = 2783.5
$ws.Range("N31").Value = -3373.5

$ws.Range("H34").Value = 2754.6274
$ws.Range("J34").Value = 2783.5
$ws.Range("L34").Value = 2783.5
$ws.Range("N34").Value = -3187.5

$ws.Range("H68").Value = 89050
$ws.Range("J68").Value = 89050
$ws.Range("L68").Value = 89050
$ws.Range("N68").Value = -90548

$ws.Range("H71").Value = 89050
$ws.Range("J71").Value = 89050
$ws.Range("L71").Value = 267150
$ws.Range("N71").Value = -274638

$ws.Range("H74").Value = 80840
$ws.Range("J74").Value = 98550
$ws.Range("L74").Value = 98550
$ws.Range("N74").Value = -100298

$ws.Range("H77").Value = 80840
$ws.Range("J77").Value = 98550
$ws.Range("L77").Value = 295650
$ws.Range("N77").Value = -304386

$ws.Range("H80").Value = 34500
$ws.Range("J80").Value = 34500
$ws.Range("L80").Value = 34500
$ws.Range("N80").Value = -36746

$ws.Range("H81").Value = 56000
$ws.Range("J81").Value = 56000
$ws.Range("L81").Value = 56000
$ws.Range("N81").Value = -57996

$ws.Range("H83").Value = 34500
$ws.Range("J83").Value = 34500
$ws.Range("L83").Value = 103500
$ws.Range("N83").Value = -114732

$ws.Range("H84").Value = 56000
$ws.Range("J84").Value = 56000
$ws.Range("L84").Value = 168000
$ws.Range("N84").Value = -177984

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 2510.8
$ws.Range("I14").Value = 2510.8
$ws.Range("K14").Value = 7532.400000000001
$ws.Range("M14").Value = -7359.400000000001

$ws.Range("H86").Value = 3668334.2
$ws.Range("J86").Value = 3668334.2
$ws.Range("L86").Value = 11005002.6
$ws.Range("N86").Value = -11007374.6

$ws.Range("H89").Value = 3668334.2
$ws.Range("J89").Value = 3668334.2
$ws.Range("L89").Value = 33015007.8
$ws.Range("N89").Value = -33026863.8

$ws.Range("H132").Value = 2439.5417
$ws.Range("I132").Value = 1018.5
$ws.Range("K132").Value = 9166.5
$ws.Range("M132").Value = -6636.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H69").Value = 39760
$ws.Range("J69").Value = 39760
$ws.Range("L69").Value = 39760
$ws.Range("N69").Value = -41258

$ws.Range("H72").Value = 39760
$ws.Range("J72").Value = 39760
$ws.Range("L72").Value = 119280
$ws.Range("N72").Value = -126768

$ws.Range("H97").Value = 4135.9
$ws.Range("I97").Value = 4787.84
$ws.Range("K97").Value = 4787.84
$ws.Range("M97").Value = -4291.84

$ws.Range("H132").Value = 265886.4
$ws.Range("I132").Value = 314920.97
$ws.Range("K132").Value = 944762.9099999999
$ws.Range("M132").Value = -942232.9099999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H9").Value = 3500

$ws.Range("H30").Value = 883
$ws.Range("I30").Value = 16
$ws.Range("K30").Value = 16
$ws.Range("M30").Value = 92

$ws.Range("H46").Value = 5356.7617
$ws.Range("I46").Value = 4346.154
$ws.Range("K46").Value = 4346.154
$ws.Range("M46").Value = -4158.154

$ws.Range("H55").Value = 522.625
$ws.Range("J55").Value = 859
$ws.Range("L55").Value = 859
$ws.Range("N55").Value = -1205

$ws.Range("H61").Value = 4266.8066
$ws.Range("I61").Value = 3040.9443
$ws.Range("J61").Value = 5964.154
$ws.Range("K61").Value = 3040.9443
$ws.Range("L61").Value = 5964.154
$ws.Range("M61").Value = -2838.9443
$ws.Range("N61").Value = -6368.154

$ws.Range("H113").Value = 4266.8066
$ws.Range("I113").Value = 3040.9443
$ws.Range("J113").Value = 5964.154
$ws.Range("K113").Value = 3040.9443
$ws.Range("L113").Value = 5964.154
$ws.Range("M113").Value = -870.9443000000001
$ws.Range("N113").Value = -10304.154

$ws.Range("H132").Value = 4127.885
$ws.Range("I132").Value = 2399.875
$ws.Range("K132").Value = 7199.625
$ws.Range("M132").Value = -4669.625

$ws.Range("H136").Value = 2921.4092
$ws.Range("I136").Value = 2254.0625
$ws.Range("K136").Value = 6762.1875
$ws.Range("M136").Value = -4212.1875

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 371.18182
$ws.Range("I100").Value = 310.3
$ws.Range("J100").Value = 980
$ws.Range("K100").Value = 620.6
$ws.Range("L100").Value = 1960
$ws.Range("M100").Value = -79.60000000000002
$ws.Range("N100").Value = -3042

$ws.Range("H107").Value = 30338.857
$ws.Range("I107").Value = 39034.52
$ws.Range("K107").Value = 117103.56
$ws.Range("M107").Value = -115183.56

$ws.Range("H113").Value = 1465.9524
$ws.Range("I113").Value = 1591.7693
$ws.Range("J113").Value = 1261.5
$ws.Range("K113").Value = 4775.3079
$ws.Range("L113").Value = 3784.5
$ws.Range("M113").Value = -2605.3079
$ws.Range("N113").Value = -8124.5
